{"js": "// Split the trailing \" Perform practical.\" run into six bold runs whose\n// concatenated text reads \" To perform this  practical click on Visualize Ckt.\"\n// (matching the target OOXML diff exactly, run-for-run).\n\nconst body = context.document.body;\nconst results = body.search(\" Perform practical.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the target text \" Perform practical.\"');\n}\n\nconst target = results.items[0];\n\n// Flat-OPC (pkg:package) snippet carrying the exact six <w:r> runs we want\n// in place of the single run that used to hold \" Perform practical.\".\n// Using insertOoxml (instead of insertText) ensures each <w:r> stays a\n// distinct run instead of being coalesced with its bold neighbours.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:rPr><w:b/></w:rPr><w:t>To p</w:t></w:r>\n            <w:r><w:rPr><w:b/></w:rPr><w:t>erform</w:t></w:r>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> this </w:t></w:r>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> practical</w:t></w:r>\n            <w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> click on Visualize Ckt.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Split the trailing \" Perform practical.\" run into six bold runs whose\n# concatenated text reads \" To perform this  practical click on Visualize Ckt.\"\n# (matching the target OOXML diff exactly, run-for-run).\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \" Perform practical.\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find the target text \" Perform practical.\"'\n}\n\n# Flat-OPC (pkg:package) snippet carrying the exact six <w:r> runs we want in\n# place of the single run that used to hold \" Perform practical.\". Using\n# InsertXML (instead of setting Range.Text) ensures each <w:r> stays a\n# distinct run instead of being coalesced with its bold neighbours.\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n          '<w:p>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t>To p</w:t></w:r>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t>erform</w:t></w:r>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> this </w:t></w:r>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> practical</w:t></w:r>' +\n            '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\"> click on Visualize Ckt.</w:t></w:r>' +\n          '</w:p>' +\n        '</w:body>' +\n      '</w:document>' +\n    '</pkg:xmlData>' +\n  '</pkg:part>' +\n'</pkg:package>'\n\n$rng.InsertXML($xml)\n"}
